$wb = $excel.ActiveWorkbook

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 21627.676
$ws.Range("I28").Value = 34550.145
$ws.Range("J28").Value = 752.9231
$ws.Range("K28").Value = 34550.145
$ws.Range("L28").Value = 752.9231
$ws.Range("M28").Value = -34065.145
$ws.Range("N28").Value = -1722.9231

# ALC row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 4328.5713
$ws.Range("J70").Value = 1800
$ws.Range("L70").Value = 5400
$ws.Range("N70").Value = -5940

# ALC row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 4328.5713
$ws.Range("J73").Value = 1800
$ws.Range("L73").Value = 5400
$ws.Range("N73").Value = -7272

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2269.9375
$ws.Range("I116").Value = 2301.111
$ws.Range("J116").Value = 2229.8572
$ws.Range("K116").Value = 2301.111
$ws.Range("L116").Value = 2229.8572
$ws.Range("M116").Value = 1140.889
$ws.Range("N116").Value = -9113.8572

# ALC row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 18503.428
$ws.Range("I125").Value = 31053
$ws.Range("J125").Value = 1770.6666
$ws.Range("K125").Value = 279477
$ws.Range("L125").Value = 15935.9994
$ws.Range("M125").Value = -277017
$ws.Range("N125").Value = -20855.9994

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2212.6365
$ws.Range("I138").Value = 1364.6818
$ws.Range("J138").Value = 3908.5454
$ws.Range("K138").Value = 4094.0454
$ws.Range("L138").Value = 11725.6362
$ws.Range("M138").Value = 1045.9546
$ws.Range("N138").Value = -22005.6362

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3174.06
$ws.Range("I141").Value = 1031.0667
$ws.Range("K141").Value = 3093.2001
$ws.Range("M141").Value = 2086.7999

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6259.2856
$ws.Range("I32").Value = 5229.255
$ws.Range("J32").Value = 16765.6
$ws.Range("K32").Value = 5229.255
$ws.Range("L32").Value = 16765.6
$ws.Range("M32").Value = -4942.255
$ws.Range("N32").Value = -17339.6

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2750
$ws.Range("I63").Value = 2750
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2750
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -2064

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2750
$ws.Range("I66").Value = 2750
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 13750
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -10318

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3745.139
$ws.Range("I74").Value = 969.8077
$ws.Range("K74").Value = 969.8077
$ws.Range("M74").Value = -95.80769999999995

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3745.139
$ws.Range("I77").Value = 969.8077
$ws.Range("K77").Value = 4849.0385
$ws.Range("M77").Value = -481.0384999999997

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 184448.17
$ws.Range("I102").Value = 334214
$ws.Range("K102").Value = 334214
$ws.Range("M102").Value = -332592

# ARM row 105
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H105").Value = 38000
$ws.Range("J105").Value = 38000
$ws.Range("L105").Value = 38000
$ws.Range("N105").Value = -44988

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1510.5
$ws.Range("I99").Value = 1010
$ws.Range("J99").Value = 2011
$ws.Range("K99").Value = 1010
$ws.Range("L99").Value = 2011
$ws.Range("M99").Value = 488
$ws.Range("N99").Value = -5007

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 50752.562
$ws.Range("I31").Value = 6057.625
$ws.Range("J31").Value = 95447.5
$ws.Range("K31").Value = 6057.625
$ws.Range("L31").Value = 95447.5
$ws.Range("M31").Value = -5762.625
$ws.Range("N31").Value = -96037.5

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 50752.562
$ws.Range("I34").Value = 6057.625
$ws.Range("J34").Value = 95447.5
$ws.Range("K34").Value = 6057.625
$ws.Range("L34").Value = 95447.5
$ws.Range("M34").Value = -5855.625
$ws.Range("N34").Value = -95851.5

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 44450736
$ws.Range("I132").Value = 74082380
$ws.Range("J132").Value = 3255.6667
$ws.Range("K132").Value = 222247140
$ws.Range("L132").Value = 9767.000100000001
$ws.Range("M132").Value = -222244610
$ws.Range("N132").Value = -14827.0001

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3107.3845
$ws.Range("I134").Value = 3199.6667
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 9599.000100000001
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -7064.000100000001
$ws.Range("N134").Value = -11070

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 627.19354
$ws.Range("J131").Value = 957.2222
$ws.Range("L131").Value = 2871.6666
$ws.Range("N131").Value = -12951.6666

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3397.1072
$ws.Range("I102").Value = 3453.25
$ws.Range("J102").Value = 3256.75
$ws.Range("K102").Value = 3453.25
$ws.Range("L102").Value = 3256.75
$ws.Range("M102").Value = -1831.25
$ws.Range("N102").Value = -6500.75

# GSM row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 277.5
$ws.Range("J107").Value = 280
$ws.Range("L107").Value = 280
$ws.Range("N107").Value = -4120

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1113.4839
$ws.Range("I7").Value = 984.9
$ws.Range("J7").Value = 1347.2727
$ws.Range("K7").Value = 984.9
$ws.Range("L7").Value = 1347.2727
$ws.Range("M7").Value = -872.9
$ws.Range("N7").Value = -1571.2727

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2126.5908
$ws.Range("I40").Value = 1763.0769
$ws.Range("J40").Value = 2651.6667
$ws.Range("K40").Value = 1763.0769
$ws.Range("L40").Value = 2651.6667
$ws.Range("M40").Value = -1627.0769
$ws.Range("N40").Value = -2923.6667

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1398.4667
$ws.Range("I93").Value = 1455.6666
$ws.Range("K93").Value = 1455.6666
$ws.Range("M93").Value = -207.6666

# LTW row 107
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H107").Value = 3500
$ws.Range("I107").Value = 3500
$ws.Range("K107").Value = 3500
$ws.Range("M107").Value = -1580

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 1113.4839
$ws.Range("I126").Value = 984.9
$ws.Range("J126").Value = 1347.2727
$ws.Range("K126").Value = 2954.7
$ws.Range("L126").Value = 4041.8181
$ws.Range("M126").Value = -484.6999999999998
$ws.Range("N126").Value = -8981.8181

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7931.05
$ws.Range("I132").Value = 10474.077
$ws.Range("J132").Value = 3208.2856
$ws.Range("K132").Value = 31422.231
$ws.Range("L132").Value = 9624.856800000001
$ws.Range("M132").Value = -28892.231
$ws.Range("N132").Value = -14684.8568

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 7765.1113
$ws.Range("I136").Value = 2398.1333
$ws.Range("J136").Value = 34600
$ws.Range("K136").Value = 7194.3999
$ws.Range("L136").Value = 103800
$ws.Range("M136").Value = -4644.3999
$ws.Range("N136").Value = -108900

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 450.92307
$ws.Range("I107").Value = 425.63635
$ws.Range("K107").Value = 1276.90905
$ws.Range("M107").Value = 643.09095

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 14287202
$ws.Range("I122").Value = 20001252
$ws.Range("J122").Value = 2075
$ws.Range("K122").Value = 60003756
$ws.Range("L122").Value = 6225
$ws.Range("M122").Value = -60001306
$ws.Range("N122").Value = -11125

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2317.675
$ws.Range("I132").Value = 2337.9285
$ws.Range("K132").Value = 7013.7855
$ws.Range("M132").Value = -4483.7855
